$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "dSF" column (F) values for a handful of rows, per repulled data.
$ws.Range("F2").Value = 9
$ws.Range("F3").Value = -11
$ws.Range("F6").Value = 6
$ws.Range("F11").Value = 9
